$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# AZ1 currently holds the "Mean" header; it becomes "Run 50" and a new
# "Mean" header is appended in BA1 (copy AZ1's formatting so BA1 matches
# the other header cells).
$ws.Range("AZ1").Copy()
$ws.Range("BA1").PasteSpecial(-4122)
$ws.Range("AZ1").Value = "Run 50"
$ws.Range("BA1").Value = "Mean"

# New "Run 50" values (previously in AZ for each row) and the
# recomputed "Mean" values (now in BA for each row).
$run50 = 175151467.5514772
$newMean = 234684594.0238228

for ($row = 2; $row -le 14; $row++) {
    $ws.Cells.Item($row, 52).Value = $run50
    $ws.Cells.Item($row, 53).Value = $newMean
}
